# Update countries & provincias Spain
# Applies the data refresh captured by the diff:
#   - Re-orders a couple of country rows (Mayotte/Eslovenia, Kenia/Sri Lanka/Letonia/Albania)
#     by writing the correct country name + stats into each row.
#   - Refreshes case/death numbers for Estados Unidos (row 4).
#   - Updates the "Datos actualizados..." timestamp in A1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Update the "last updated" timestamp banner ---
$ws.Range("A1").Value = "Datos actualizados a 20 de Mayo de 2020 a las 15:35"

# --- Row 4: Estados Unidos ---
$ws.Range("B4").Value = 1572091
$ws.Range("C4").Value = 1508
$ws.Range("D4").Value = 361227
$ws.Range("E4").Value = 1117270
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 61
$ws.Range("H4").Value = 93594

# --- Row 98: becomes Mayotte (new data) ---
$ws.Range("A98").Value = "Mayotte"
$ws.Range("B98").Value = 1475
$ws.Range("C98").Value = 105
$ws.Range("D98").Value = 894
$ws.Range("E98").Value = 562
$ws.Range("F98").Value = 0
$ws.Range("G98").Value = 1
$ws.Range("H98").Value = 19

# --- Row 99: becomes Eslovenia (previous row 98 data) ---
$ws.Range("A99").Value = "Eslovenia"
$ws.Range("B99").Value = 1468
$ws.Range("C99").Value = 1
$ws.Range("D99").Value = 1340
$ws.Range("E99").Value = 23
$ws.Range("F99").Value = 0
$ws.Range("G99").Value = 1
$ws.Range("H99").Value = 105

# --- Row 105: becomes Kenia (new data) ---
$ws.Range("A105").Value = "Kenia"
$ws.Range("B105").Value = 1029
$ws.Range("C105").Value = 66
$ws.Range("D105").Value = 358
$ws.Range("E105").Value = 621
$ws.Range("F105").Value = 0
$ws.Range("G105").Value = 0
$ws.Range("H105").Value = 50

# --- Row 106: becomes Sri Lanka (previous row 105 data) ---
$ws.Range("A106").Value = "Sri Lanka"
$ws.Range("B106").Value = 1027
$ws.Range("C106").Value = 4
$ws.Range("D106").Value = 584
$ws.Range("E106").Value = 434
$ws.Range("F106").Value = 0
$ws.Range("G106").Value = 0
$ws.Range("H106").Value = 9

# --- Row 107: becomes Letonia (previous row 106 data) ---
$ws.Range("A107").Value = "Letonia"
$ws.Range("B107").Value = 1016
$ws.Range("C107").Value = 4
$ws.Range("D107").Value = 694
$ws.Range("E107").Value = 301
$ws.Range("F107").Value = 0
$ws.Range("G107").Value = 0
$ws.Range("H107").Value = 21

# --- Row 108: becomes Albania (previous row 107 data) ---
$ws.Range("A108").Value = "Albania"
$ws.Range("B108").Value = 964
$ws.Range("C108").Value = 15
$ws.Range("D108").Value = 758
$ws.Range("E108").Value = 175
$ws.Range("F108").Value = 0
$ws.Range("G108").Value = 0
$ws.Range("H108").Value = 31
